# Power Rankings update script
$wb = $excel.ActiveWorkbook
$points = $wb.Worksheets.Item("POINTS")

# --- Enter WK8 (column J) scores for rows 3-16 ---
$points.Range("J3").Value = 93.7
$points.Range("J4").Value = 81.5
$points.Range("J5").Value = 106.7
$points.Range("J6").Value = 119.2
$points.Range("J7").Value = 117.6
$points.Range("J8").Value = 118.2
$points.Range("J9").Value = 93.4
$points.Range("J10").Value = 94.3
$points.Range("J11").Value = 99
$points.Range("J12").Value = 92.8
$points.Range("J13").Value = 97.7
$points.Range("J14").Value = 115.5
$points.Range("J15").Value = 82.7
$points.Range("J16").Value = 107.2

# --- Update total-wins (column C) for rows 22-35 ---
$points.Range("C22").Value = 7
$points.Range("C24").Value = 2
$points.Range("C25").Value = 4
$points.Range("C27").Value = 5
$points.Range("C29").Value = 2
$points.Range("C33").Value = 6
$points.Range("C35").Value = 6

# --- Clear roster scores (column G) for rows 22-35 ---
$points.Range("G22:G35").ClearContents()
